# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '90.289.30'
$ws.Range("E2").Value = '  -1.01%  '

# Row 3
$ws.Range("D3").Value = '3.085.88'
$ws.Range("E3").Value = '  -2.08%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.06'
$ws.Range("E5").Value = '  +2.77%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '614.39'
$ws.Range("E6").Value = '  -1.68%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.05'
$ws.Range("E7").Value = '  +15.32%  '

# Row 8
$ws.Range("E8").Value = '  -7.17%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.00%  '

# Row 10
$ws.Range("D10").Value = '3.081.19'
$ws.Range("E10").Value = '  -2.07%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.712'
$ws.Range("E11").Value = '  -2.31%  '

# Row 12
$ws.Range("E12").Value = '  +2.73%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000240'
$ws.Range("E13").Value = '  -5.75%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.44'
$ws.Range("E14").Value = '  -0.31%  '

# Row 15
$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.50'
$ws.Range("E15").Value = '  +1.20%  '

# Row 16
$ws.Range("D16").Value = '90.036.99'
$ws.Range("E16").Value = '  -1.38%  '

# Row 17
$ws.Range("D17").Value = '3.643.89'
$ws.Range("E17").Value = '  -2.04%  '

# Row 18
$ws.Range("D18").Value = '3.070.86'
$ws.Range("E18").Value = '  -2.21%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.63'
$ws.Range("E19").Value = '  -6.44%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.13'
$ws.Range("E20").Value = '  -1.04%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000204'
$ws.Range("E21").Value = '  -10.92%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '445.12'
$ws.Range("E22").Value = '  +2.54%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.87'
$ws.Range("E23").Value = '  +0.81%  '

# Row 24
$ws.Range("E24").Value = '  +4.01%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.98'
$ws.Range("E25").Value = '  -1.97%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.27'
$ws.Range("E26").Value = '  +7.78%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.82'
$ws.Range("E27").Value = '  -5.41%  '

# Row 28
$ws.Range("D28").Value = '3.232.30'
$ws.Range("E28").Value = '  -1.83%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.15%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.24'
$ws.Range("E30").Value = '  +2.32%  '

# Row 31
$ws.Range("E31").Value = '  -6.04%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.63'
$ws.Range("E32").Value = '  +17.09%  '

# Row 33
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.939'
$ws.Range("E33").Value = '  -5.96%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.198'
$ws.Range("E34").Value = '  +36.10%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.147'
$ws.Range("E35").Value = '  +4.66%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '499.33'
$ws.Range("E36").Value = '  -6.78%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.89'
$ws.Range("E37").Value = '  +1.13%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.88'
$ws.Range("E38").Value = '  -6.15%  '

# Row 39
$ws.Range("E39").Value = '  -2.65%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.48'
$ws.Range("E40").Value = '  -11.79%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.423'
$ws.Range("E41").Value = '  +11.59%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.21'
$ws.Range("E42").Value = '  -0.41%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0852'
$ws.Range("E43").Value = '  +8.16%  '

# Row 44
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("E45").Value = '  +0.05%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.01'
$ws.Range("E46").Value = '  +20.07%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.691'
$ws.Range("E47").Value = '  +10.41%  '

# Row 48
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '148.62'
$ws.Range("E48").Value = '  +3.31%  '

# Row 49
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.57'
$ws.Range("E49").Value = '  +9.24%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.54'
$ws.Range("E50").Value = '  +0.77%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.33'
$ws.Range("E51").Value = '  +2.17%  '
